# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
#
# The worksheet lists one debt row per "Periodo Mora" (YYMM code) for a
# single worker, from the oldest period (1908) at row 16 down to the
# newest (2207) at row 51. This edit flips that list so the newest
# period is on top and the oldest is at the bottom (new statement
# periods added on top, older ones pushed down toward the end), while
# keeping each period's own "Valor Mora" (column F) value attached to
# it as it moves. Every period carries 32648 except the newest one,
# 2207, which carries 27207 - so after the flip that value now sits at
# the top (row 16) instead of the bottom (row 51).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Periods exactly as originally listed top-to-bottom (row 16 .. row 51).
$periods = @(
    "1908","1909","1910","1911","1912",
    "2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112",
    "2201","2202","2203","2204","2205","2206","2207"
)

$startRow = 16
$count = $periods.Length

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $srcIndex = $count - 1 - $i
    $period = $periods[$srcIndex]

    if ($period -eq "2207") {
        $value = 27207
    } else {
        $value = 32648
    }

    $ws.Cells.Item($row, 5).Value = $period
    $ws.Cells.Item($row, 6).Value = $value
}
